# Updated symbol list on Sat Dec 31 09:57:16 UTC 2022 with GitHub Actions
# Re-applies the refreshed Coinranking snapshot values (prices re-sorted / re-scraped)
# onto the existing worksheet rows, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'245.50"

# Row 3
$ws.Range("D3").Value = "'26.33"

# Row 4
$ws.Range("D4").Value = "'5.138"

# Row 5
$ws.Range("D5").Value = "'0.05588"

# Row 6
$ws.Range("D6").Value = "'6.490"

# Row 8
$ws.Range("D8").Value = "'0.8160"

# Row 9
$ws.Range("D9").Value = "'0.8480"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1344"
$ws.Range("E10").Value = "9WazirXWRX"

# Row 11
$ws.Range("B11").Value = "One"
$ws.Range("C11").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D11").Value = "'0.003530"
$ws.Range("E11").Value = "10OneONEBestin24h"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.06953"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

# Row 13
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03219"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02851"
$ws.Range("E14").Value = "13BitrueCoinBTR"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09403"
$ws.Range("E15").Value = "14BitMartTokenBMX"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001518"
$ws.Range("E16").Value = "15BitForexTokenBF"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006245"
$ws.Range("E17").Value = "16TigerCashTCH"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.548"
$ws.Range("E18").Value = "17LEOLEO"

# Row 19
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.118"
$ws.Range("E19").Value = "18BTSETokenBTSE"

# Row 20
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3134"
$ws.Range("E20").Value = "19BitpandaEcosystemTokenBEST"

# Row 22
$ws.Range("D22").Value = "'3.757"

# Row 25
$ws.Range("D25").Value = "'0.001249"

# Row 26
$ws.Range("D26").Value = "'0.004606"

# Row 27
$ws.Range("D27").Value = "'0.00009602"

# Row 40
$ws.Range("D40").Value = "'0.03653"

# Row 41
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1361"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.006132"
$ws.Range("E42").Value = "41KickTokenKICK"

# Row 43
$ws.Range("D43").Value = "'0.002558"

# Row 44
$ws.Range("D44").Value = "'0.007953"

# Row 45
$ws.Range("D45").Value = "'0.00005317"

# Row 49
$ws.Range("D49").Value = "'0.00002100"

# Row 50
$ws.Range("D50").Value = "'0.0002000"
